$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 17 de Julio de 2020 a las 20:03
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 17 de Julio de 2020 a las 20:03'

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 3725956
$ws.Cells.Item(4, 3).Value = 30931
$ws.Cells.Item(4, 4).Value = 1694496
$ws.Cells.Item(4, 5).Value = 1889942
$ws.Cells.Item(4, 7).Value = 400
$ws.Cells.Item(4, 8).Value = 141518

# Row 6: India
$ws.Cells.Item(6, 2).Value = 1039045
$ws.Cells.Item(6, 3).Value = 33408
$ws.Cells.Item(6, 4).Value = 654056
$ws.Cells.Item(6, 5).Value = 358704
$ws.Cells.Item(6, 7).Value = 676
$ws.Cells.Item(6, 8).Value = 26285

# Row 18: Turquia
$ws.Cells.Item(18, 2).Value = 217799
$ws.Cells.Item(18, 3).Value = 926
$ws.Cells.Item(18, 4).Value = 199834
$ws.Cells.Item(18, 5).Value = 12507
$ws.Cells.Item(18, 7).Value = 18
$ws.Cells.Item(18, 8).Value = 5458

# Row 19: Alemania
$ws.Cells.Item(19, 2).Value = 202025
$ws.Cells.Item(19, 3).Value = 189
$ws.Cells.Item(19, 5).Value = 5966
$ws.Cells.Item(19, 7).Value = 2
$ws.Cells.Item(19, 8).Value = 9159

# Row 58: Irlanda
$ws.Cells.Item(58, 2).Value = 25730
$ws.Cells.Item(58, 3).Value = 32
$ws.Cells.Item(58, 5).Value = 614
$ws.Cells.Item(58, 7).Value = 3
$ws.Cells.Item(58, 8).Value = 1752

# Row 65: Marruecos
$ws.Cells.Item(65, 2).Value = 16726
$ws.Cells.Item(65, 3).Value = 181
$ws.Cells.Item(65, 4).Value = 14360
$ws.Cells.Item(65, 5).Value = 2102
$ws.Cells.Item(65, 7).Value = 1
$ws.Cells.Item(65, 8).Value = 264

# Row 67: Uzbekistan
$ws.Cells.Item(67, 2).Value = 15482
$ws.Cells.Item(67, 3).Value = 416
$ws.Cells.Item(67, 5).Value = 6622
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(67, 8).Value = 77

# Row 108: Maldivas
$ws.Cells.Item(108, 2).Value = 2913
$ws.Cells.Item(108, 3).Value = 14
$ws.Cells.Item(108, 4).Value = 2340
$ws.Cells.Item(108, 5).Value = 558

# Row 113: Sri Lanka
$ws.Cells.Item(113, 2).Value = 2697
$ws.Cells.Item(113, 3).Value = 10
$ws.Cells.Item(113, 5).Value = 674

# Row 114: Congo
$ws.Cells.Item(114, 1).Value = 'Congo'
$ws.Cells.Item(114, 2).Value = 2633
$ws.Cells.Item(114, 3).Value = 275
$ws.Cells.Item(114, 4).Value = 626
$ws.Cells.Item(114, 5).Value = 1958
$ws.Cells.Item(114, 7).Value = 1
$ws.Cells.Item(114, 8).Value = 49

# Row 115: Mali
$ws.Cells.Item(115, 1).Value = 'Mali'
$ws.Cells.Item(115, 2).Value = 2467
$ws.Cells.Item(115, 3).Value = 27
$ws.Cells.Item(115, 4).Value = 1791
$ws.Cells.Item(115, 5).Value = 555
$ws.Cells.Item(115, 8).Value = 121

# Row 116: Cuba
$ws.Cells.Item(116, 1).Value = 'Cuba'
$ws.Cells.Item(116, 2).Value = 2444
$ws.Cells.Item(116, 3).Value = 4
$ws.Cells.Item(116, 4).Value = 2300
$ws.Cells.Item(116, 5).Value = 57
$ws.Cells.Item(116, 8).Value = 87

# Row 126: Hong Kong
$ws.Cells.Item(126, 5).Value = 439
$ws.Cells.Item(126, 7).Value = 1
$ws.Cells.Item(126, 8).Value = 11

# Row 127: Libia
$ws.Cells.Item(127, 1).Value = 'Libia'
$ws.Cells.Item(127, 2).Value = 1704
$ws.Cells.Item(127, 3).Value = 52
$ws.Cells.Item(127, 4).Value = 380
$ws.Cells.Item(127, 5).Value = 1277
$ws.Cells.Item(127, 8).Value = 47

# Row 128: Sierra Leona
$ws.Cells.Item(128, 1).Value = 'Sierra Leona'
$ws.Cells.Item(128, 2).Value = 1688
$ws.Cells.Item(128, 3).Value = 10
$ws.Cells.Item(128, 4).Value = 1219
$ws.Cells.Item(128, 5).Value = 404
$ws.Cells.Item(128, 7).Value = 1
$ws.Cells.Item(128, 8).Value = 65

# Row 129: Suazilandia
$ws.Cells.Item(129, 1).Value = 'Suazilandia'
$ws.Cells.Item(129, 2).Value = 1619
$ws.Cells.Item(129, 3).Value = 67
$ws.Cells.Item(129, 4).Value = 769
$ws.Cells.Item(129, 5).Value = 829
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 21

# Row 130: Benin
$ws.Cells.Item(130, 1).Value = 'Benin'
$ws.Cells.Item(130, 2).Value = 1602
$ws.Cells.Item(130, 3).Value = 139
$ws.Cells.Item(130, 4).Value = 782
$ws.Cells.Item(130, 5).Value = 789
$ws.Cells.Item(130, 7).Value = 3
$ws.Cells.Item(130, 8).Value = 31

# Row 131: Yemen
$ws.Cells.Item(131, 2).Value = 1576
$ws.Cells.Item(131, 3).Value = 24
$ws.Cells.Item(131, 5).Value = 441
$ws.Cells.Item(131, 7).Value = 2
$ws.Cells.Item(131, 8).Value = 440

# Row 141: Namibia
$ws.Cells.Item(141, 1).Value = 'Namibia'
$ws.Cells.Item(141, 2).Value = 1078
$ws.Cells.Item(141, 3).Value = 46
$ws.Cells.Item(141, 4).Value = 32
$ws.Cells.Item(141, 5).Value = 1044
$ws.Cells.Item(141, 8).Value = 2

# Row 142: Uganda
$ws.Cells.Item(142, 1).Value = 'Uganda'
$ws.Cells.Item(142, 2).Value = 1056
$ws.Cells.Item(142, 3).Value = 5
$ws.Cells.Item(142, 4).Value = 1023
$ws.Cells.Item(142, 5).Value = 33
$ws.Cells.Item(142, 8).Value = 0

# Row 143: Burkina Faso
$ws.Cells.Item(143, 1).Value = 'Burkina Faso'
$ws.Cells.Item(143, 2).Value = 1045
$ws.Cells.Item(143, 3).Value = 7
$ws.Cells.Item(143, 4).Value = 887
$ws.Cells.Item(143, 5).Value = 105
$ws.Cells.Item(143, 8).Value = 53

# Row 144: Republica de Chipre
$ws.Cells.Item(144, 2).Value = 1033
$ws.Cells.Item(144, 3).Value = 2
$ws.Cells.Item(144, 5).Value = 169

# Row 174: Guadalupe
$ws.Cells.Item(174, 2).Value = 195
$ws.Cells.Item(174, 3).Value = 5
$ws.Cells.Item(174, 4).Value = 172
$ws.Cells.Item(174, 5).Value = 9

# Row 209: Islas Malvinas
$ws.Cells.Item(209, 1).Value = 'Islas Malvinas'

# Row 210: Groenlandia
$ws.Cells.Item(210, 1).Value = 'Groenlandia'
